$wb = $excel.ActiveWorkbook

# --- PRESENTATION sheet (sheet16): add the new task rows ---
$presentation = $wb.Worksheets.Item("PRESENTATION")

$rows = @(
    "Create an outline of what we are presenting that matches the web site",
    "Create a visual slide show using some technology ",
    "Each slide should visually display our `"brand`" and be focused on the topic ",
    "Be willing to flex and change as needed based on the needs of the project",
    "Be sure to include all the visualizations possible.",
    "Be sure to include the topic of the machine learning model",
    "Be sure to include conclusions of the project",
    "Be sure to proof read and provide credits for any quotes taken",
    "Share with the team ",
    "publish to the git. "
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $presentation.Cells.Item($r, 1).Value = $i + 1
    $presentation.Cells.Item($r, 2).Value = $rows[$i]
}

# --- DEFINE THE QUESTIONS sheet (sheet10): change the selection ---
$questions = $wb.Worksheets.Item("DEFINE THE QUESTIONS")
$questions.Activate() | Out-Null
$questions.Range("A1:C5").Select() | Out-Null

# --- PRESENTATION sheet: leave the selection parked on K18 ---
$presentation.Activate() | Out-Null
$presentation.Range("K18").Select() | Out-Null

# restore the originally active sheet/tab so it stays marked as selected
$wb.Worksheets.Item("TASKS").Activate() | Out-Null
